# Generate Report for Handoff
# The "bab6b1a6-5c71-44a9-bc79-338ec00640e3.md" file moved from "In Translation"
# to "Ready for handoff" in both the zh-cn and de-de localization sheets, with
# a new Priority ("mt") and refreshed "Latest Handoff Datetime" / Latest HO
# Xliff Generate Date timestamps. Reflect the same status + date refresh on
# the Overview summary sheet.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row for bab6b1a6-...md is row 3 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("H3").Value = "2016-08-17 06:11:13"

# --- de-de sheet: row for bab6b1a6-...md is row 3 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("H3").Value = "2016-08-17 06:11:17"

# --- Overview sheet: row for bab6b1a6-...md is row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 06:11:17"

# Column widths auto-fit to the wider "Ready for handoff" text.
$wsOverview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$wsOverview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$wsZh.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$wsDe.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
